# Weekly update: a new price observation is inserted at the top of the
# "Feria Lagunitas de Puerto Montt - Perejil" log (row 217), pushing the
# previous rows 217-219 down to 218-220.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 217 - this shifts the existing rows
# 217, 218, 219 down to 218, 219, 220 (and grows the used range to R220).
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with this week's data.
$ws.Cells.Item(217, 1).Value  = 4
$ws.Cells.Item(217, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(217, 3).Value  = "Los Lagos"
$ws.Cells.Item(217, 4).Value  = 44628
$ws.Cells.Item(217, 5).Value  = 10
$ws.Cells.Item(217, 6).Value  = 100112044
$ws.Cells.Item(217, 7).Value  = "Perejil"
$ws.Cells.Item(217, 8).Value  = "Sin especificar"
$ws.Cells.Item(217, 9).Value  = "Primera"
$ws.Cells.Item(217, 10).Value = 180
$ws.Cells.Item(217, 11).Value = 6000
$ws.Cells.Item(217, 12).Value = 6000
$ws.Cells.Item(217, 13).Value = 6000
$ws.Cells.Item(217, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(217, 15).Value = "Región Metropolitana"
$ws.Cells.Item(217, 16).Value = 2000
$ws.Cells.Item(217, 17).Value = 3
$ws.Cells.Item(217, 18).Value = "Hortaliza"
